# Scheduled-runner update: refresh cached Market Board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# profitability sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 10001350
$ws.Range("I43").Value = 12501250
$ws.Range("K43").Value = 12501250
$ws.Range("M43").Value = -12501181

$ws.Range("H53").Value = 55555828
$ws.Range("J53").Value = 111111500
$ws.Range("L53").Value = 111111500
$ws.Range("N53").Value = -111112774

$ws.Range("H58").Value = 569.6
$ws.Range("I58").Value = 612.25
$ws.Range("J58").Value = 399
$ws.Range("K58").Value = 1836.75
$ws.Range("L58").Value = 1197
$ws.Range("M58").Value = -1686.75
$ws.Range("N58").Value = -1497

$ws.Range("H106").Value = 1244.9429
$ws.Range("I106").Value = 1325.8518
$ws.Range("J106").Value = 971.875
$ws.Range("K106").Value = 1325.8518
$ws.Range("L106").Value = 971.875
$ws.Range("M106").Value = -694.8517999999999
$ws.Range("N106").Value = -2233.875

$ws.Range("H107").Value = 13889199
$ws.Range("I107").Value = 333.1
$ws.Range("K107").Value = 333.1
$ws.Range("M107").Value = 1586.9

$ws.Range("H116").Value = 5211.143
$ws.Range("I116").Value = 5211.143
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5211.143
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1769.143
$ws.Range("N116").ClearContents()

$ws.Range("H138").Value = 5286.17
$ws.Range("I138").Value = 26475
$ws.Range("J138").Value = 3556.4695
$ws.Range("K138").Value = 79425
$ws.Range("L138").Value = 10669.4085
$ws.Range("M138").Value = -74285
$ws.Range("N138").Value = -20949.4085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 202926.6
$ws.Range("I32").Value = 239349.62
$ws.Range("K32").Value = 239349.62
$ws.Range("M32").Value = -239062.62

$ws.Range("H61").Value = 1524580.5
$ws.Range("I61").Value = 48030.305
$ws.Range("J61").Value = 3522266
$ws.Range("K61").Value = 48030.305
$ws.Range("L61").Value = 3522266
$ws.Range("M61").Value = -47818.305
$ws.Range("N61").Value = -3522690

$ws.Range("H74").Value = 497727.88
$ws.Range("I74").Value = 1134.7307
$ws.Range("J74").Value = 1177276.4
$ws.Range("K74").Value = 1134.7307
$ws.Range("L74").Value = 1177276.4
$ws.Range("M74").Value = -260.7307000000001
$ws.Range("N74").Value = -1179024.4

$ws.Range("H77").Value = 497727.88
$ws.Range("I77").Value = 1134.7307
$ws.Range("J77").Value = 1177276.4
$ws.Range("K77").Value = 5673.6535
$ws.Range("L77").Value = 5886382
$ws.Range("M77").Value = -1305.6535
$ws.Range("N77").Value = -5895118

$ws.Range("H122").Value = 1000.7647
$ws.Range("I122").Value = 969.5625
$ws.Range("K122").Value = 2908.6875
$ws.Range("M122").Value = -458.6875

$ws.Range("H136").Value = 1524580.5
$ws.Range("I136").Value = 48030.305
$ws.Range("J136").Value = 3522266
$ws.Range("K136").Value = 144090.915
$ws.Range("L136").Value = 10566798
$ws.Range("M136").Value = -141540.915
$ws.Range("N136").Value = -10571898

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 89998.7
$ws.Range("J20").Value = 89998.7
$ws.Range("L20").Value = 89998.7
$ws.Range("N20").Value = -90470.7

$ws.Range("H30").Value = 89998.7
$ws.Range("J30").Value = 89998.7
$ws.Range("L30").Value = 89998.7
$ws.Range("N30").Value = -90180.7

$ws.Range("H31").Value = 3425.3125
$ws.Range("I31").Value = 3602.1
$ws.Range("J31").Value = 3344.9546
$ws.Range("K31").Value = 3602.1
$ws.Range("L31").Value = 3344.9546
$ws.Range("M31").Value = -3307.1
$ws.Range("N31").Value = -3934.9546

$ws.Range("H34").Value = 3425.3125
$ws.Range("I34").Value = 3602.1
$ws.Range("J34").Value = 3344.9546
$ws.Range("K34").Value = 3602.1
$ws.Range("L34").Value = 3344.9546
$ws.Range("M34").Value = -3400.1
$ws.Range("N34").Value = -3748.9546

$ws.Range("H58").Value = 2078.12
$ws.Range("I58").Value = 2435.2856
$ws.Range("J58").Value = 1623.5454
$ws.Range("K58").Value = 2435.2856
$ws.Range("L58").Value = 1623.5454
$ws.Range("M58").Value = -2232.2856
$ws.Range("N58").Value = -2029.5454

$ws.Range("H128").Value = 89998.7
$ws.Range("J128").Value = 89998.7
$ws.Range("L128").Value = 89998.7
$ws.Range("N128").Value = -99958.7

$ws.Range("H134").Value = 1435.037
$ws.Range("I134").Value = 1347.5769
$ws.Range("K134").Value = 4042.7307
$ws.Range("M134").Value = -1507.7307

$ws.Range("H136").Value = 2078.12
$ws.Range("I136").Value = 2435.2856
$ws.Range("J136").Value = 1623.5454
$ws.Range("K136").Value = 7305.8568
$ws.Range("L136").Value = 4870.6362
$ws.Range("M136").Value = -4755.8568
$ws.Range("N136").Value = -9970.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 15.8
$ws.Range("I23").Value = 25.5
$ws.Range("J23").Value = 9.333333
$ws.Range("K23").Value = 76.5
$ws.Range("L23").Value = 27.999999
$ws.Range("M23").Value = 158.5
$ws.Range("N23").Value = -497.999999

$ws.Range("H50").Value = 6725
$ws.Range("I50").Value = 3963.1667
$ws.Range("K50").Value = 11889.5001
$ws.Range("M50").Value = -11408.5001

$ws.Range("H53").Value = 6725
$ws.Range("I53").Value = 3963.1667
$ws.Range("K53").Value = 11889.5001
$ws.Range("M53").Value = -11408.5001

$ws.Range("H118").Value = 5410.5454
$ws.Range("I118").Value = 6417.778
$ws.Range("K118").Value = 19253.334
$ws.Range("M118").Value = -18010.334

$ws.Range("H119").Value = 111127000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2249.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2249.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2249.5
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2481.5

$ws.Range("H36").Value = 2917.1667
$ws.Range("J36").Value = 2700.75
$ws.Range("L36").Value = 2700.75
$ws.Range("N36").Value = -3670.75

$ws.Range("H80").Value = 10152592
$ws.Range("J80").Value = 35840772
$ws.Range("L80").Value = 35840772
$ws.Range("N80").Value = -35842768

$ws.Range("H83").Value = 10152592
$ws.Range("J83").Value = 35840772
$ws.Range("L83").Value = 179203860
$ws.Range("N83").Value = -179213844

$ws.Range("H93").Value = 81288.25
$ws.Range("J93").Value = 68384.336
$ws.Range("L93").Value = 68384.336
$ws.Range("N93").Value = -72128.336

$ws.Range("H113").Value = 1890.8
$ws.Range("I113").Value = 1838.625
$ws.Range("J113").Value = 2099.5
$ws.Range("K113").Value = 1838.625
$ws.Range("L113").Value = 2099.5
$ws.Range("M113").Value = 331.375
$ws.Range("N113").Value = -6439.5

$ws.Range("H126").Value = 1877.25
$ws.Range("I126").Value = 2048.5
$ws.Range("K126").Value = 6145.5
$ws.Range("M126").Value = -3675.5

$ws.Range("H132").Value = 695173.2
$ws.Range("J132").Value = 2139516.5
$ws.Range("L132").Value = 6418549.5
$ws.Range("N132").Value = -6423609.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 54899.5

$ws.Range("H46").Value = 13683.923
$ws.Range("I46").Value = 21248.5
$ws.Range("J46").Value = 7200
$ws.Range("K46").Value = 21248.5
$ws.Range("L46").Value = 7200
$ws.Range("M46").Value = -21060.5
$ws.Range("N46").Value = -7576

$ws.Range("H55").Value = 1347.025
$ws.Range("I55").Value = 1367.9565
$ws.Range("K55").Value = 1367.9565
$ws.Range("M55").Value = -1194.9565

$ws.Range("H68").Value = 3592
$ws.Range("J68").Value = 4148.5
$ws.Range("L68").Value = 4148.5
$ws.Range("N68").Value = -5646.5

$ws.Range("H71").Value = 3592
$ws.Range("J71").Value = 4148.5
$ws.Range("L71").Value = 20742.5
$ws.Range("N71").Value = -28230.5

$ws.Range("H136").Value = 1637.2043
$ws.Range("I136").Value = 1576.6383
$ws.Range("J136").Value = 1699.0869
$ws.Range("K136").Value = 4729.9149
$ws.Range("L136").Value = 5097.2607
$ws.Range("M136").Value = -2179.9149
$ws.Range("N136").Value = -10197.2607

$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 50000
$ws.Range("K137").Value = 50000
$ws.Range("M137").Value = -44900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 627.7273
$ws.Range("I113").Value = 627.7273
$ws.Range("K113").Value = 1883.1819
$ws.Range("M113").Value = 286.8181

$ws.Range("H119").Value = 50345
$ws.Range("J119").Value = 50345
$ws.Range("L119").Value = 50345
$ws.Range("N119").Value = -60021

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
